# ------------------------------------------------------------------
# Edit: rotate three worksheet names (customers/purchases/vendors)
# and swap their contents so that:
#   - sheet at position 1 (was "customers") becomes "purchases" and
#     gets filled with purchase-order data
#   - sheet at position 4 (was "vendors") becomes "customers" and
#     gets a fresh (header-only) customers table
#   - sheet at position 5 (was "purchases") becomes "vendors" and
#     gets the old vendors header row (moved, not retyped)
# sheets 2 (inventory), 3 (products), 6 (sales) are untouched.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets (position-based, collision-safe order) ------
# customers(1) -> purchases, vendors(4) -> customers, purchases(5) -> vendors
# is a 3-cycle, so free up a name with a scratch name first.
$wb.Worksheets.Item(4).Name = "__scratch_rename__"
$wb.Worksheets.Item(5).Name = "vendors"
$wb.Worksheets.Item(1).Name = "purchases"
$wb.Worksheets.Item(4).Name = "customers"

$wsPurchases = $wb.Worksheets.Item(1)
$wsCustomers = $wb.Worksheets.Item(4)
$wsVendors   = $wb.Worksheets.Item(5)

# --- 2. sheet5 (now "vendors"): reuse the old vendors header row --
# Before sheet4 ("vendors") is overwritten below, copy its header
# row (A1:F1) onto sheet5, which still holds the old purchases
# header (A1:G1) that is no longer wanted.
$wsCustomers.Range("A1:F1").Copy()
$wsVendors.Range("A1").PasteSpecial(-4104)
$excel.CutCopyMode = $false
# drop the leftover 7th column ("Total_Cost") from the old purchases header
$wsVendors.Range("G1").Clear()

# --- 3. sheet4 (now "customers"): header-only customers table -----
$wsCustomers.Range("A1").Value = "Customer_ID"
$wsCustomers.Range("B1").Value = "Name"
$wsCustomers.Range("C1").Value = "Email"
$wsCustomers.Range("D1").Value = "Address"
# remove the old vendors columns E1:F1 which are no longer part of
# this (narrower) table
$wsCustomers.Range("E1:F1").Clear()

# --- 4. sheet1 (now "purchases"): purchases data -------------------
$wsPurchases.Range("A1").Value = "HSN_Code"
$wsPurchases.Range("B1").Value = "Product_Name"
$wsPurchases.Range("C1").Value = "Vendor"
$wsPurchases.Range("D1").Value = "Date"
$wsPurchases.Range("E1").Value = "Units"
$wsPurchases.Range("F1").Value = "Cost_per_Unit"
$wsPurchases.Range("G1").Value = "Total_Cost"
# give the new header cells (E1:G1) the same bold/boxed style as the
# rest of row 1, reusing the existing header style (cellXfs index 1)
$wsPurchases.Range("A1").Copy()
$wsPurchases.Range("E1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 -- every value typed in as literal text (matches source data
# entry, including the numeric-looking ones and the date)
$wsPurchases.Range("A2:G2").NumberFormat = "@"
$wsPurchases.Range("A2").Value = "1"
$wsPurchases.Range("B2").Value = "Keyboard"
$wsPurchases.Range("C2").Value = "Amit"
$wsPurchases.Range("D2").Value = "2025-06-08"
$wsPurchases.Range("E2").Value = "100"
$wsPurchases.Range("F2").Value = "1200"
$wsPurchases.Range("G2").Value = "120000"
$wsPurchases.Range("A2:G2").Style = "Normal"

# Row 3 -- numbers stored as real numbers, text stays text
$wsPurchases.Range("B3:D3").NumberFormat = "@"
$wsPurchases.Range("A3").Value = 2
$wsPurchases.Range("B3").Value = "mouse"
$wsPurchases.Range("C3").Value = "Amit"
$wsPurchases.Range("D3").Value = "2025-06-08"
$wsPurchases.Range("E3").Value = 21
$wsPurchases.Range("F3").Value = 2200
$wsPurchases.Range("G3").Value = 46200
$wsPurchases.Range("B3:D3").Style = "Normal"

Write-Host "Sheets now: $(($wb.Worksheets | ForEach-Object { $_.Name }) -join ', ')"
